# Adds a new test-result row (44) to the "fullDetails" sheet: a "signIn" test
# run that failed, copied from the structure of the previous row, plus its
# screenshot hyperlink - this mirrors the author's fix for Windows/Linux
# screenshot-link compatibility (new row uses a Windows-style "C:\..." path).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 43
$dstRow = 44

# 1) Clone the previous row's cells verbatim (values only, via copy/paste-
#    special) so every text/number/boolean cell keeps its original storage
#    type instead of being re-inferred by a plain .Value assignment.
$ws.Range("A$srcRow`:AT$srcRow").Copy()
$ws.Range("A$dstRow").PasteSpecial(-4163)

# Helper: write $value into $cellRef while forcing it to stay a literal text
# cell (never a number/bool), regardless of what it looks like - builds it as
# a quoted-string formula result in a scratch cell, off in an unused column,
# then copies that result in as a value and wipes the scratch cell again.
function Set-TextValue($cellRef, $value) {
    $scratch = $ws.Range("ZZ100")
    $escaped = $value.Replace("""", """""")
    $scratch.Formula = "=""$escaped"""
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# 2) Overwrite the handful of cells that actually differ for this new run.
Set-TextValue "T44"  "signIn, beton.yatsuk@gmail.com, ABC123, Hi, Avner!"
Set-TextValue "Z44"  "signIn"
Set-TextValue "AG44" "2015-11-30:12-45-23"
Set-TextValue "AJ44" "1448887523167"
Set-TextValue "AM44" "2015-11-30 14:46:01"
Set-TextValue "AO44" "fail"
Set-TextValue "AT44" "C:\Users\AvnerG\git\Beton\Beton/test-output/screenshots2015-11-30-14-46-01-199-IST.png"

# 3) Re-add the hyperlink on the new screenshot cell (Windows-style path).
#    Hyperlinks.Add auto-applies the blue/underlined "Hyperlink" cell style;
#    the sheet's existing hyperlink cells (e.g. AT43) carry no such styling,
#    so strip it back off to match.
$ws.Hyperlinks.Add($ws.Range("AT44"), "C:\Users\AvnerG\git\Beton\Beton/test-output/screenshots2015-11-30-14-46-01-199-IST.png")
$ws.Range("AT44").ClearFormats()
